# Scheduled-runner refresh of market-price-derived columns (H:N) on the
# Ragnarok profit sheets. Source item/recipe columns (A:G) are untouched;
# only the computed current-price / profit columns are refreshed per row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 531.25
$ws.Range("I33").Value = 425.33334
$ws.Range("J33").Value = 849
$ws.Range("K33").Value = 425.33334
$ws.Range("L33").Value = 849
$ws.Range("M33").Value = -196.33334
$ws.Range("N33").Value = -1307

$ws.Range("H112").Value = 4432.3335
$ws.Range("J112").Value = 4960.385
$ws.Range("L112").Value = 14881.155
$ws.Range("N112").Value = -17097.155

$ws.Range("H115").Value = 1912.2307
$ws.Range("I115").Value = 478.66666
$ws.Range("K115").Value = 1435.99998
$ws.Range("M115").Value = 131.0000199999999

$ws.Range("H118").Value = 1285.7778
$ws.Range("I118").Value = 912.8333
$ws.Range("J118").Value = 2031.6666
$ws.Range("K118").Value = 2738.4999
$ws.Range("L118").Value = 6094.9998
$ws.Range("M118").Value = -1081.4999
$ws.Range("N118").Value = -9408.9998

$ws.Range("H138").Value = 4084.0833
$ws.Range("I138").Value = 2038.5
$ws.Range("J138").Value = 6947.9
$ws.Range("K138").Value = 6115.5
$ws.Range("L138").Value = 20843.7
$ws.Range("M138").Value = -975.5
$ws.Range("N138").Value = -31123.7

$ws.Range("H141").Value = 7944.2
$ws.Range("I141").Value = 8049.222
$ws.Range("K141").Value = 24147.666
$ws.Range("M141").Value = -18967.666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 788.2105
$ws.Range("I2").Value = 864.4286
$ws.Range("K2").Value = 864.4286
$ws.Range("M2").Value = -751.4286

$ws.Range("H61").Value = 7247899
$ws.Range("I61").Value = 8340232
$ws.Range("K61").Value = 8340232
$ws.Range("M61").Value = -8340020

$ws.Range("H74").Value = 1762.2963
$ws.Range("I74").Value = 1620.0834
$ws.Range("K74").Value = 1620.0834
$ws.Range("M74").Value = -746.0834

$ws.Range("H77").Value = 1762.2963
$ws.Range("I77").Value = 1620.0834
$ws.Range("K77").Value = 8100.416999999999
$ws.Range("M77").Value = -3732.416999999999

$ws.Range("H116").Value = 788.2105
$ws.Range("I116").Value = 864.4286
$ws.Range("K116").Value = 864.4286
$ws.Range("M116").Value = 1429.5714

$ws.Range("H132").Value = 5266283
$ws.Range("I132").Value = 3019.7334
$ws.Range("J132").Value = 25003520
$ws.Range("K132").Value = 9059.200199999999
$ws.Range("L132").Value = 75010560
$ws.Range("M132").Value = -6529.200199999999
$ws.Range("N132").Value = -75015620

$ws.Range("H136").Value = 7247899
$ws.Range("I136").Value = 8340232
$ws.Range("K136").Value = 25020696
$ws.Range("M136").Value = -25018146

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 788.2105
$ws.Range("I3").Value = 864.4286
$ws.Range("K3").Value = 864.4286
$ws.Range("M3").Value = -750.4286

$ws.Range("H99").Value = 2081.1428
$ws.Range("J99").Value = 2971.4285
$ws.Range("L99").Value = 2971.4285
$ws.Range("N99").Value = -5967.4285

$ws.Range("H134").Value = 3846853.8
$ws.Range("I134").Value = 562.375
$ws.Range("K134").Value = 1687.125
$ws.Range("M134").Value = 847.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 43482188
$ws.Range("I31").Value = 111114450
$ws.Range("J31").Value = 4304.4287
$ws.Range("K31").Value = 111114450
$ws.Range("L31").Value = 4304.4287
$ws.Range("M31").Value = -111114155
$ws.Range("N31").Value = -4894.4287

$ws.Range("H34").Value = 43482188
$ws.Range("I34").Value = 111114450
$ws.Range("J34").Value = 4304.4287
$ws.Range("K34").Value = 111114450
$ws.Range("L34").Value = 4304.4287
$ws.Range("M34").Value = -111114248
$ws.Range("N34").Value = -4708.4287

$ws.Range("H62").Value = 18532448
$ws.Range("I62").Value = 7349.8335
$ws.Range("K62").Value = 7349.8335
$ws.Range("M62").Value = -6725.8335

$ws.Range("H65").Value = 18532448
$ws.Range("I65").Value = 7349.8335
$ws.Range("K65").Value = 36749.1675
$ws.Range("M65").Value = -33629.1675

$ws.Range("H99").Value = 12981.37
$ws.Range("I99").Value = 7360.5
$ws.Range("J99").Value = 21157.182
$ws.Range("K99").Value = 7360.5
$ws.Range("L99").Value = 21157.182
$ws.Range("M99").Value = -5862.5
$ws.Range("N99").Value = -24153.182

$ws.Range("H107").Value = 2414.7273
$ws.Range("I107").Value = 916.4
$ws.Range("K107").Value = 916.4
$ws.Range("M107").Value = 1003.6

$ws.Range("H126").Value = 12981.37
$ws.Range("I126").Value = 7360.5
$ws.Range("J126").Value = 21157.182
$ws.Range("K126").Value = 22081.5
$ws.Range("L126").Value = 63471.546
$ws.Range("M126").Value = -19611.5
$ws.Range("N126").Value = -68411.546

$ws.Range("H132").Value = 2723.7
$ws.Range("I132").Value = 2755.9473
$ws.Range("J132").Value = 2111
$ws.Range("K132").Value = 8267.841899999999
$ws.Range("L132").Value = 6333
$ws.Range("M132").Value = -5737.841899999999
$ws.Range("N132").Value = -11393

$ws.Range("H134").Value = 3169.6
$ws.Range("I134").Value = 2683.2727
$ws.Range("K134").Value = 8049.8181
$ws.Range("M134").Value = -5514.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 815.8095
$ws.Range("J113").Value = 680.5714
$ws.Range("L113").Value = 2041.7142
$ws.Range("N113").Value = -6381.7142

$ws.Range("H120").Value = 29434
$ws.Range("I120").Value = 27484.5
$ws.Range("K120").Value = 82453.5
$ws.Range("M120").Value = -77615.5

$ws.Range("H129").Value = 6413572
$ws.Range("I129").Value = 12502367
$ws.Range("J129").Value = 4314.579
$ws.Range("K129").Value = 37507101
$ws.Range("L129").Value = 12943.737
$ws.Range("M129").Value = -37502101
$ws.Range("N129").Value = -22943.737

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 15500
$ws.Range("I10").Value = 1001
$ws.Range("J10").Value = 29999
$ws.Range("K10").Value = 1001
$ws.Range("L10").Value = 29999
$ws.Range("M10").Value = -832
$ws.Range("N10").Value = -30337

$ws.Range("H80").Value = 3998.5
$ws.Range("I80").Value = 3998
$ws.Range("K80").Value = 3998
$ws.Range("M80").Value = -3000

$ws.Range("H83").Value = 3998.5
$ws.Range("I83").Value = 3998
$ws.Range("K83").Value = 19990
$ws.Range("M83").Value = -14998

$ws.Range("H102").Value = 3986.5
$ws.Range("I102").Value = 3942.5
$ws.Range("K102").Value = 3942.5
$ws.Range("M102").Value = -2320.5

$ws.Range("H107").Value = 246.57143
$ws.Range("I107").Value = 246.57143
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 246.57143
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1673.42857
$ws.Range("N107").ClearContents()

$ws.Range("H126").Value = 5009.346
$ws.Range("J126").Value = 5658.8
$ws.Range("L126").Value = 16976.4
$ws.Range("N126").Value = -21916.4

$ws.Range("H132").Value = 5560247.5
$ws.Range("I132").Value = 4932.643
$ws.Range("K132").Value = 14797.929
$ws.Range("M132").Value = -12267.929

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7549.892
$ws.Range("I7").Value = 7431.3477
$ws.Range("K7").Value = 7431.3477
$ws.Range("M7").Value = -7319.3477

$ws.Range("H46").Value = 942.8125
$ws.Range("J46").Value = 1999.3334
$ws.Range("L46").Value = 1999.3334
$ws.Range("N46").Value = -2375.3334

$ws.Range("H55").Value = 1425.1
$ws.Range("I55").Value = 1154.5
$ws.Range("J55").Value = 1695.7
$ws.Range("K55").Value = 1154.5
$ws.Range("L55").Value = 1695.7
$ws.Range("M55").Value = -981.5
$ws.Range("N55").Value = -2041.7

$ws.Range("H126").Value = 7549.892
$ws.Range("I126").Value = 7431.3477
$ws.Range("K126").Value = 22294.0431
$ws.Range("M126").Value = -19824.0431

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

$ws.Range("H96").Value = 7081.077
$ws.Range("I96").Value = 6972
$ws.Range("J96").Value = 7255.6
$ws.Range("K96").Value = 6972
$ws.Range("L96").Value = 7255.6
$ws.Range("M96").Value = -5599
$ws.Range("N96").Value = -10001.6
